# "Base de datos Colombia" - column header cleanups + sheet rename
$wb = $excel.ActiveWorkbook

# Inflacion sheet: "Meta inflación" -> "Meta" (table header, column B)
$wsInflacion = $wb.Worksheets.Item("Inflacion")
$wsInflacion.Range("B1").Value = "Meta"

# PIB sheet: "Delta PIB" -> "Delta_PIB" (table header, column C) to match the
# underscore naming convention used by the other delta columns
$wsPIB = $wb.Worksheets.Item("PIB")
$wsPIB.Range("C1").Value = "Delta_PIB"

# Rename "Tasas_de_interes" sheet to the singular "Tasa_de_interes"
$wsTasas = $wb.Worksheets.Item("Tasas_de_interes")
$wsTasas.Name = "Tasa_de_interes"

# Make the PIB sheet the active tab/selection (matches the saved workbook view)
$wsPIB.Activate() | Out-Null
$wsPIB.Range("C1").Select() | Out-Null
